# Applies the quiz-02 edits described by the commit diff:
#  1. Q3 ("...of the action potential?"): drop the trailing bold space
#     run and leave a collapsed "_GoBack" bookmark in its place.
#  2. Q6: merge the three runs (split by now-removed proofErr markers)
#     into a single run with text "...except ????."
#  3. Q7: merge the runs (and drop the old "_GoBack" bookmark that used
#     to sit here) into a single run "...a/an ???? receptor."
#  4. The "presynaptic" list item: drop the gramStart/gramEnd proofErr
#     wrapper around the lone run.
#  5. Q11: merge the runs into a single run ending in "...inside-out?"
#  6. Footer page-number field cached result: "4" -> "1".

$d = $word.ActiveDocument

# --- 1. Q3: trailing " " run -> collapsed "_GoBack" bookmark ----------
$rng = $d.Content
$null = $rng.Find.Execute("of the action potential?", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rng)
$rng.Text = ""

$p18 = $d.Paragraphs(18)
$pr18 = $p18.Range
$trailing = $d.Range($pr18.End - 2, $pr18.End - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}

# --- 2. Q6: merge runs into one, proofErr markers vanish with them ----
$q6old = "6. All of the following monoamine neurotransmitters are " + `
    "released from nuclei located in the midbrain and brainstem, " + `
    "except ?" + "???."
$q6new = "6. All of the following monoamine neurotransmitters are " + `
    "released from nuclei located in the midbrain and brainstem, " + `
    "except ????."
$rng = $d.Content
$null = $rng.Find.Execute($q6old, $true, $false, $false, $false, $false, `
    $true, 1, $false, $q6new, 2)

# --- 3. Q7: merge runs, implicitly removing the old "_GoBack" bookmark
$q7old = "7. With the exception of a single serotonin receptor type, " + `
    "all of the monoamines bind to a/" + "an ?" + "???" + " receptor."
$q7new = "7. With the exception of a single serotonin receptor type, " + `
    "all of the monoamines bind to a/an ???? receptor."
$rng = $d.Content
$null = $rng.Find.Execute($q7old, $true, $false, $false, $false, $false, `
    $true, 1, $false, $q7new, 2)

# --- 4. "presynaptic": drop the surrounding proofErr markers ----------
$peOld = "voltage-gated" + [char]13 + "presynaptic"
$rng = $d.Content
$null = $rng.Find.Execute($peOld, $true, $false, $false, $false, $false, `
    $true, 1, $false, $peOld, 2)

# --- 5. Q11: merge runs into one -------------------------------------
$q11old = "11. Why does Gilmore say the retina is physiologically " + `
    "backwards and anatomically " + "inside-out" + "?"
$rng = $d.Content
$null = $rng.Find.Execute($q11old, $true, $false, $false, $false, $false, `
    $true, 1, $false, $q11old, 2)

# --- 6. Footer page-number field cached result: 4 -> 1 ----------------
$sec = $d.Sections.Item(1)
$footers = $sec.Footers
$footer = $footers.Item(1)
$flds = $footer.Range.Fields
$fld = $flds.Item(1)
$fres = $fld.Result
$null = $fres.Find.Execute("4", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1", 2)

Write-Output "done"
